$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 7.62x39 damage increase: the "ammo_k_hit" coefficient (column H) for the
# two 7.62x39 rounds (rows 19 = FMJ/Perf, 20 = AP) goes from 1.05 to 1.17.
# Dependent formulas (E, J, K columns) recalc automatically.
$ws.Range("H19").Value = 1.17
$ws.Range("H20").Value = 1.17

# Update the saved view state to match (scroll position / selection).
$ws.Range("K13").Select()
